# Applies the "Add files via upload" commit: rename the sheet and refresh
# the repayment figures for the 2025-09-01..2025-09-17 cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the worksheet tab (matches the workbook.xml <sheet name=.../> diff) ---
$ws.Name = "repayment_20250901_20250917 (1)"

# --- Row 2 : Debora Retima Sihombing ---
$ws.Range("H2").Value = 17.84

# --- Row 3 : Romli ---
$ws.Range("H3").Value = 24.265999999999998

# --- Row 4 : Aldi Taufik ---
$ws.Range("D4").Value = 46
$ws.Range("E4").Value = "24,923,826.00"
$ws.Range("G4").Value = "6.83"
$ws.Range("H4").Value = 13.693

# --- Row 5 : Yandi Nugraha ---
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = "30,156,872.00"
$ws.Range("G5").Value = "9.26"
$ws.Range("H5").Value = 21.375
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = "5.57"
$ws.Range("L5").Value = "5.98"

# --- Row 6 : Axl Wicaksono ---
$ws.Range("H6").Value = 15.552

# --- Row 7 : Riska Nurlita ---
$ws.Range("D7").Value = 66
$ws.Range("E7").Value = "31,887,424.00"
$ws.Range("G7").Value = "8.59"
$ws.Range("H7").Value = 12.106

# --- Row 8 : Annisa Putri Restu ---
$ws.Range("H8").Value = 24.422999999999998
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = "3.87"
$ws.Range("L8").Value = "5.40"

# --- Row 9 : Azizah Rahmawati ---
$ws.Range("D9").Value = 39
$ws.Range("E9").Value = "36,439,115.00"
$ws.Range("G9").Value = "9.53"
$ws.Range("H9").Value = 12.930999999999999
$ws.Range("J9").Value = 9
$ws.Range("K9").Value = "5.05"
$ws.Range("L9").Value = "3.16"

# --- Row 10 : Erlangga Hutama ---
$ws.Range("D10").Value = 26
$ws.Range("E10").Value = "17,483,715.00"
$ws.Range("G10").Value = "5.70"
$ws.Range("H10").Value = 13.705
$ws.Range("J10").Value = 7
$ws.Range("K10").Value = "3.40"
$ws.Range("L10").Value = "3.21"

# --- Row 11 : Erick Ervan Dewanggga ---
$ws.Range("D11").Value = 47
$ws.Range("E11").Value = "42,116,879.00"
$ws.Range("G11").Value = "11.92"
$ws.Range("H11").Value = 13.712

# --- Row 12 : Ridhoi Berkat Zebua ---
$ws.Range("H12").Value = 19.82

# --- Row 13 : Fadilah Damayanti ---
$ws.Range("H13").Value = 17.297999999999998

# --- Row 14 : Nur Halim ---
$ws.Range("H14").Value = 10.535

# --- Row 15 : Adistira Winditya P ---
$ws.Range("D15").Value = 40
$ws.Range("E15").Value = "28,316,229.00"
$ws.Range("G15").Value = "8.00"
$ws.Range("H15").Value = 11.661
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = "3.83"
$ws.Range("L15").Value = "4.24"

# --- Row 16 : Sucika Wardani ---
$ws.Range("H16").Value = 10.638999999999999

# --- Row 17 : Wasti Feronika Sihombing ---
$ws.Range("D17").Value = 45
$ws.Range("E17").Value = "32,411,629.00"
$ws.Range("G17").Value = "9.50"
$ws.Range("H17").Value = 20.623000000000001

# --- Row 18 : Nuraini ---
$ws.Range("H18").Value = 11.128
